# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1 (08:03 -> 09:03)
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 09:03"

# Row 4 - Estados Unidos (USA)
$ws.Range("B4").Value = 1212955
$ws.Range("C4").Value = 120
$ws.Range("E4").Value = 954962
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 69925

# Row 38 - Ucrania
$ws.Range("B38").Value = 12697
$ws.Range("C38").Value = 366
$ws.Range("D38").Value = 1875
$ws.Range("E38").Value = 10506
$ws.Range("F38").Value = 167
$ws.Range("G38").Value = 13
$ws.Range("H38").Value = 316

# Row 48 - Chequia
$ws.Range("B48").Value = 7841
$ws.Range("C48").Value = 22
$ws.Range("D48").Value = 3816
$ws.Range("E48").Value = 3773
$ws.Range("F48").Value = 59

# Row 60 - Kazajistan
$ws.Range("B60").Value = 4121
$ws.Range("C60").Value = 72
$ws.Range("E60").Value = 2919

# Row 111 - Georgia
$ws.Range("B111").Value = 604
$ws.Range("C111").Value = 11
$ws.Range("D111").Value = 240
$ws.Range("E111").Value = 355

# Row 146 - Brunei
$ws.Range("D146").Value = 131
$ws.Range("E146").Value = 6
